$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $range = $ws.Range($CellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $Value
    $range.Style = $origStyle
}

Set-TextValue "D2" '23.319.41'
Set-TextValue "E2" '  -0.54%  '
Set-TextValue "D3" '1.624.44'
Set-TextValue "E3" '  -0.30%  '
Set-TextValue "E4" '  +0.56%  '
Set-TextValue "E5" '  +0.50%  '
Set-TextValue "D6" '303.81'
Set-TextValue "E6" '  -1.21%  '
Set-TextValue "D7" '0.3778'
Set-TextValue "E7" '  -0.07%  '
Set-TextValue "D8" '51.85'
Set-TextValue "E8" '  -2.38%  '
Set-TextValue "D9" '0.3605'
Set-TextValue "E9" '  -1.72%  '
Set-TextValue "E10" '  -4.31%  '
Set-TextValue "D11" '0.08071'
Set-TextValue "E11" '  -1.57%  '
Set-TextValue "E12" '  +0.54%  '
Set-TextValue "E13" '  -3.22%  '
Set-TextValue "D14" '6.540'
Set-TextValue "E14" '  -2.04%  '
Set-TextValue "D15" '0.00001242'
Set-TextValue "E15" '  -1.80%  '
Set-TextValue "D16" '7.204'
Set-TextValue "E16" '  -3.54%  '
Set-TextValue "D17" '1.624.62'
Set-TextValue "E17" '  -0.30%  '
Set-TextValue "D18" '93.29'
Set-TextValue "E18" '  -1.66%  '
Set-TextValue "D19" '0.06908'
Set-TextValue "E19" '  -0.36%  '
Set-TextValue "D20" '17.91'
Set-TextValue "E20" '  -2.87%  '
Set-TextValue "D21" '1.002'
Set-TextValue "D22" '6.431'
Set-TextValue "E22" '  -2.43%  '
Set-TextValue "D23" '23.346.87'
Set-TextValue "E23" '  -0.50%  '
Set-TextValue "E24" '  -2.20%  '
Set-TextValue "E25" '  +2.39%  '
Set-TextValue "D26" '2.443'
Set-TextValue "E26" '  +0.24%  '
Set-TextValue "D27" '21.04'
Set-TextValue "E27" '  -1.75%  '
Set-TextValue "D28" '148.86'
Set-TextValue "E28" '  -1.33%  '
Set-TextValue "D29" '5.288'
Set-TextValue "E29" '  -0.50%  '
Set-TextValue "D30" '134.43'
Set-TextValue "E30" '  -1.45%  '
Set-TextValue "D31" '2.293'
Set-TextValue "E31" '  -5.73%  '
Set-TextValue "D32" '1.805.78'
Set-TextValue "E32" '  -0.04%  '
Set-TextValue "D33" '6.726'
Set-TextValue "E33" '  -3.18%  '
Set-TextValue "D34" '10.82'
Set-TextValue "E34" '  +3.17%  '
Set-TextValue "E35" '  -3.53%  '
Set-TextValue "D36" '0.02811'
Set-TextValue "E36" '  +0.11%  '
Set-TextValue "E37" '  -0.67%  '
Set-TextValue "D38" '0.08807'
Set-TextValue "E38" '  -0.46%  '
Set-TextValue "D39" '6.091'
Set-TextValue "E39" '  -2.63%  '
Set-TextValue "D40" '0.07083'
Set-TextValue "E40" '  -5.64%  '
Set-TextValue "D41" '1.358'
Set-TextValue "E41" '  -3.65%  '
Set-TextValue "D42" '0.7026'
Set-TextValue "E42" '  -2.07%  '
Set-TextValue "D43" '16.14'
Set-TextValue "E43" '  -0.35%  '
Set-TextValue "D44" '12.22'
Set-TextValue "E44" '  -4.66%  '
Set-TextValue "E45" '  -3.09%  '
Set-TextValue "D46" '1.000'
Set-TextValue "E46" '  +0.48%  '
Set-TextValue "D47" '2.310'
Set-TextValue "E47" '  -2.26%  '
Set-TextValue "D48" '3.981'
Set-TextValue "E48" '  -1.33%  '
Set-TextValue "D49" '0.07967'
Set-TextValue "E49" '  -0.76%  '
Set-TextValue "D50" '1.198'
Set-TextValue "E50" '  -1.49%  '
Set-TextValue "D51" '125.53'
Set-TextValue "E51" '  -5.01%  '
